$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26/27 swap: Monero <-> PancakeSwap
# Row 38/39 swap: ImmutableX <-> TrustWalletToken
$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.14"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.30%  "
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "167.48"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.26%  "
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.46"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.42%  "
$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.884"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.42%  "

# Price / Volume(1h) updates for remaining rows
$ws.Range("D2").Value = "36.937.51"
$ws.Range("E2").Value = "  +4.46%  "
$ws.Range("D3").Value = "1.914.13"
$ws.Range("E3").Value = "  +1.66%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "248.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.24%  "
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "47.71"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +10.47%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.377"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +6.64%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "58.42"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.48%  "
$ws.Range("E11").Value = "  +2.61%  "
$ws.Range("E12").Value = "  +2.19%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.36"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +12.40%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.820"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.77%  "
$ws.Range("D15").Value = "2.194.39"
$ws.Range("E15").Value = "  +1.70%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.12"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.39%  "
$ws.Range("D17").Value = "1.914.50"
$ws.Range("E17").Value = "  +1.24%  "
$ws.Range("D18").Value = "36.908.85"
$ws.Range("E18").Value = "  +4.50%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "74.56"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.77%  "
$ws.Range("D20").Value = "0.0₃0856"
$ws.Range("E20").Value = "  +3.98%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.59"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +6.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "251.31"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.15"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.46"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -6.09%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.81"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.68%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.66"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.37%  "
$ws.Range("E30").Value = "  +1.00%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.55"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.36%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0609"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.28%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0907"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +26.26%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.28"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.36%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.90"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.27%  "
$ws.Range("E36").Value = "  +0.10%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.95"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +38.86%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.97"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.96%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "104.59"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.48%  "
$ws.Range("E42").Value = "  +3.53%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.48"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.50%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.87"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +19.30%  "
$ws.Range("E45").Value = "  +2.33%  "
$ws.Range("D46").Value = "1.349.45"
$ws.Range("E46").Value = "  +2.74%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.42"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.80%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0835"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.04%  "
$ws.Range("E49").Value = "  +3.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.39"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.03%  "
$ws.Range("E51").Value = "  +13.77%  "
